$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.066.93'
$ws.Range('E2').Value = '  -1.27%  '
$ws.Range('D3').Value = '3.150.18'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.94'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.66'
$ws.Range('E6').Value = '  -3.10%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.145.40'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.39'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.466'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  -2.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.02'
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').Value = '3.671.51'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').Value = '64.113.79'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '3.141.21'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '489.05'
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.70'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.712'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.64'
$ws.Range('E23').Value = '  -5.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.60'
$ws.Range('E24').Value = '  +3.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.32'
$ws.Range('E25').Value = '  -3.39%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.77'
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.21'
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.00'
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.39'
$ws.Range('E31').Value = '  +3.10%  '
$ws.Range('E32').Value = '  -6.29%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.03'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.70'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  -4.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.97'
$ws.Range('E39').Value = '  -8.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '437.22'
$ws.Range('E40').Value = '  -5.22%  '
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.30'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('D44').Value = '2.924.81'
$ws.Range('E44').Value = '  +2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.259'
$ws.Range('E45').Value = '  -3.94%  '
$ws.Range('E46').Value = '  -6.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.97'
$ws.Range('E49').Value = '  -2.76%  '
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.24'
$ws.Range('E51').Value = '  -0.44%  '
